$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.658.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "'3.577.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'242.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").Value = "'653.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'1.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.64%  "
$ws.Range("D8").Value = "'0.412"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("E9").Value = "  +8.27%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "'3.577.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'43.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'6.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "'4.242.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'96.482.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "'0.0000257"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "'3.583.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'7.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'12.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "'18.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "'0.537"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.48%  "
$ws.Range("D23").Value = "'508.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'6.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.80%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").Value = "'96.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").Value = "'13.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.36%  "
$ws.Range("D29").Value = "'3.769.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'3.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.152"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.57%  "
$ws.Range("D32").Value = "'11.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.98%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'31.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'8.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.13%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'626.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.56%  "
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +12.19%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.151"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "'0.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'1.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.87%  "
$ws.Range("D45").Value = "'5.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.17%  "
$ws.Range("D46").Value = "'0.0433"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("D47").Value = "'2.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("D48").Value = "'23.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'32.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.83%  "
$ws.Range("D50").Value = "'3.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'8.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.87%  "
